$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 618
    5  = 12868
    8  = 508
    9  = 472
    10 = 1164
    11 = 966
    12 = 13712
    13 = 14180
    22 = 1079
    25 = 935
    26 = 5269
    27 = 10
    28 = 282
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
